# Regenerate merged AHB files
# 1) Rename the "_old"/"_new" suffixed header labels to "_FV2310"/"_FV2404"
# 2) Turn the header row + data range into an Excel Table (ListObject)
# 3) Freeze the header row (pane split below row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$suffixMap = @{
    "Segmentname_old"         = "Segmentname_FV2310"
    "Segmentgruppe_old"       = "Segmentgruppe_FV2310"
    "Segment_old"             = "Segment_FV2310"
    "Datenelement_old"        = "Datenelement_FV2310"
    "Segment ID_old"          = "Segment ID_FV2310"
    "Code_old"                = "Code_FV2310"
    "Qualifier_old"           = "Qualifier_FV2310"
    "Beschreibung_old"        = "Beschreibung_FV2310"
    "Bedingungsausdruck_old"  = "Bedingungsausdruck_FV2310"
    "Bedingung_old"           = "Bedingung_FV2310"
    "Segmentname_new"         = "Segmentname_FV2404"
    "Segmentgruppe_new"       = "Segmentgruppe_FV2404"
    "Segment_new"             = "Segment_FV2404"
    "Datenelement_new"        = "Datenelement_FV2404"
    "Segment ID_new"          = "Segment ID_FV2404"
    "Code_new"                = "Code_FV2404"
    "Qualifier_new"           = "Qualifier_FV2404"
    "Beschreibung_new"        = "Beschreibung_FV2404"
    "Bedingungsausdruck_new"  = "Bedingungsausdruck_FV2404"
    "Bedingung_new"           = "Bedingung_FV2404"
}

for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $old = $cell.Value2
    if ($suffixMap.ContainsKey($old)) {
        $cell.Value2 = $suffixMap[$old]
    }
}

# Convert the used range into a native Excel Table with an autofilter
$tableRange = $ws.Range("A1:U86")
$listObject = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$listObject.Name = "Table1"

# Freeze the header row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
